$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (year 2025) values per the latest data refresh
$ws.Range("B7").Value = 2142476.95
$ws.Range("C7").Value = -52.57313037997906
$ws.Range("D7").Value = 2014
$ws.Range("E7").Value = 2014
$ws.Range("F7").Value = 1063.791931479643
$ws.Range("G7").Value = 9.783548246543017
